$d = $word.ActiveDocument
$d.Content.Find.Execute("83-48=35", $true, $false, $false, $false, $false, $true, 1, $false, "4+8=12", 1) | Out-Null
$d.Content.Find.Execute("85-44=41", $true, $false, $false, $false, $false, $true, 1, $false, "81-39=42", 1) | Out-Null
$d.Content.Find.Execute("4+22=26", $true, $false, $false, $false, $false, $true, 1, $false, "62-60=2", 1) | Out-Null
$d.Content.Find.Execute("76-74=2", $true, $false, $false, $false, $false, $true, 1, $false, "91-5=86", 1) | Out-Null
$d.Content.Find.Execute("58-24=34", $true, $false, $false, $false, $false, $true, 1, $false, "12+31=43", 1) | Out-Null
$d.Content.Find.Execute("46+48=94", $true, $false, $false, $false, $false, $true, 1, $false, "2+77=79", 1) | Out-Null
$d.Content.Find.Execute("46+22=68", $true, $false, $false, $false, $false, $true, 1, $false, "68-30=38", 1) | Out-Null
$d.Content.Find.Execute("31-26=5", $true, $false, $false, $false, $false, $true, 1, $false, "52+43=95", 1) | Out-Null
$d.Content.Find.Execute("5+61=66", $true, $false, $false, $false, $false, $true, 1, $false, "75-13=62", 1) | Out-Null
$d.Content.Find.Execute("93-11=82", $true, $false, $false, $false, $false, $true, 1, $false, "33+62=95", 1) | Out-Null
$d.Content.Find.Execute("64-51=13", $true, $false, $false, $false, $false, $true, 1, $false, "50+11=61", 1) | Out-Null
$d.Content.Find.Execute("56+42=98", $true, $false, $false, $false, $false, $true, 1, $false, "38+37=75", 1) | Out-Null
$d.Content.Find.Execute("72-35=37", $true, $false, $false, $false, $false, $true, 1, $false, "37-36=1", 1) | Out-Null
$d.Content.Find.Execute("5+30=35", $true, $false, $false, $false, $false, $true, 1, $false, "60-27=33", 1) | Out-Null
$d.Content.Find.Execute("97-25=72", $true, $false, $false, $false, $false, $true, 1, $false, "2+5=7", 1) | Out-Null
$d.Content.Find.Execute("10+84=94", $true, $false, $false, $false, $false, $true, 1, $false, "93+3=96", 1) | Out-Null
$d.Content.Find.Execute("31+43=74", $true, $false, $false, $false, $false, $true, 1, $false, "67+8=75", 1) | Out-Null
$d.Content.Find.Execute("94-84=10", $true, $false, $false, $false, $false, $true, 1, $false, "34+55=89", 1) | Out-Null
$d.Content.Find.Execute("64-53=11", $true, $false, $false, $false, $false, $true, 1, $false, "83-78=5", 1) | Out-Null
$d.Content.Find.Execute("20-19=1", $true, $false, $false, $false, $false, $true, 1, $false, "14+15=29", 1) | Out-Null
$d.Content.Find.Execute("46-38=8", $true, $false, $false, $false, $false, $true, 1, $false, "64-48=16", 1) | Out-Null
$d.Content.Find.Execute("96-13=83", $true, $false, $false, $false, $false, $true, 1, $false, "93-54=39", 1) | Out-Null
$d.Content.Find.Execute("82+1=83", $true, $false, $false, $false, $false, $true, 1, $false, "61-14=47", 1) | Out-Null
$d.Content.Find.Execute("13+0=13", $true, $false, $false, $false, $false, $true, 1, $false, "80-36=44", 1) | Out-Null
$d.Content.Find.Execute("27-24=3", $true, $false, $false, $false, $false, $true, 1, $false, "78+9=87", 1) | Out-Null
$d.Content.Find.Execute("50+15=65", $true, $false, $false, $false, $false, $true, 1, $false, "69-69=0", 1) | Out-Null
$d.Content.Find.Execute("0+55=55", $true, $false, $false, $false, $false, $true, 1, $false, "48+7=55", 1) | Out-Null
$d.Content.Find.Execute("1+39=40", $true, $false, $false, $false, $false, $true, 1, $false, "0+93=93", 1) | Out-Null
$d.Content.Find.Execute("32+55=87", $true, $false, $false, $false, $false, $true, 1, $false, "15+45=60", 1) | Out-Null
$d.Content.Find.Execute("62-30=32", $true, $false, $false, $false, $false, $true, 1, $false, "15+57=72", 1) | Out-Null
$d.Content.Find.Execute("74+11=85", $true, $false, $false, $false, $false, $true, 1, $false, "92-9=83", 1) | Out-Null
$d.Content.Find.Execute("45-31=14", $true, $false, $false, $false, $false, $true, 1, $false, "36+50=86", 1) | Out-Null
$d.Content.Find.Execute("72-43=29", $true, $false, $false, $false, $false, $true, 1, $false, "18+54=72", 1) | Out-Null
$d.Content.Find.Execute("58+0=58", $true, $false, $false, $false, $false, $true, 1, $false, "58+37=95", 1) | Out-Null
$d.Content.Find.Execute("31-2=29", $true, $false, $false, $false, $false, $true, 1, $false, "25+35=60", 1) | Out-Null
$d.Content.Find.Execute("86-39=47", $true, $false, $false, $false, $false, $true, 1, $false, "36+32=68", 1) | Out-Null
$d.Content.Find.Execute("98-79=19", $true, $false, $false, $false, $false, $true, 1, $false, "91-54=37", 1) | Out-Null
$d.Content.Find.Execute("44-30=14", $true, $false, $false, $false, $false, $true, 1, $false, "18+57=75", 1) | Out-Null
$d.Content.Find.Execute("89-66=23", $true, $false, $false, $false, $false, $true, 1, $false, "51+10=61", 1) | Out-Null
$d.Content.Find.Execute("76-53=23", $true, $false, $false, $false, $false, $true, 1, $false, "80-13=67", 1) | Out-Null
$d.Content.Find.Execute("99-55=44", $true, $false, $false, $false, $false, $true, 1, $false, "75-17=58", 1) | Out-Null
$d.Content.Find.Execute("17+74=91", $true, $false, $false, $false, $false, $true, 1, $false, "84-15=69", 1) | Out-Null
$d.Content.Find.Execute("67-47=20", $true, $false, $false, $false, $false, $true, 1, $false, "6+83=89", 1) | Out-Null
$d.Content.Find.Execute("97-62=35", $true, $false, $false, $false, $false, $true, 1, $false, "5-1=4", 1) | Out-Null
$d.Content.Find.Execute("54-27=27", $true, $false, $false, $false, $false, $true, 1, $false, "86-18=68", 1) | Out-Null
$d.Content.Find.Execute("3+76=79", $true, $false, $false, $false, $false, $true, 1, $false, "17+57=74", 1) | Out-Null
$d.Content.Find.Execute("98-71=27", $true, $false, $false, $false, $false, $true, 1, $false, "62-40=22", 1) | Out-Null
$d.Content.Find.Execute("22+40=62", $true, $false, $false, $false, $false, $true, 1, $false, "30+47=77", 1) | Out-Null
$d.Content.Find.Execute("99-3=96", $true, $false, $false, $false, $false, $true, 1, $false, "71+1=72", 1) | Out-Null
$d.Content.Find.Execute("97-42=55", $true, $false, $false, $false, $false, $true, 1, $false, "49+15=64", 1) | Out-Null
$d.Content.Find.Execute("34-12=22", $true, $false, $false, $false, $false, $true, 1, $false, "64+10=74", 1) | Out-Null
$d.Content.Find.Execute("91-44=47", $true, $false, $false, $false, $false, $true, 1, $false, "72-71=1", 1) | Out-Null
$d.Content.Find.Execute("18+17=35", $true, $false, $false, $false, $false, $true, 1, $false, "46+2=48", 1) | Out-Null
$d.Content.Find.Execute("40+49=89", $true, $false, $false, $false, $false, $true, 1, $false, "7+9=16", 1) | Out-Null
$d.Content.Find.Execute("78-58=20", $true, $false, $false, $false, $false, $true, 1, $false, "17+78=95", 1) | Out-Null
$d.Content.Find.Execute("37+7=44", $true, $false, $false, $false, $false, $true, 1, $false, "18+44=62", 1) | Out-Null
$d.Content.Find.Execute("85-38=47", $true, $false, $false, $false, $false, $true, 1, $false, "40-18=22", 1) | Out-Null
$d.Content.Find.Execute("66-22=44", $true, $false, $false, $false, $false, $true, 1, $false, "97-37=60", 1) | Out-Null
$d.Content.Find.Execute("72+6=78", $true, $false, $false, $false, $false, $true, 1, $false, "56+38=94", 1) | Out-Null
$d.Content.Find.Execute("32+34=66", $true, $false, $false, $false, $false, $true, 1, $false, "3+80=83", 1) | Out-Null
$d.Content.Find.Execute("97-88=9", $true, $false, $false, $false, $false, $true, 1, $false, "18+22=40", 1) | Out-Null
$d.Content.Find.Execute("14+24=38", $true, $false, $false, $false, $false, $true, 1, $false, "8+73=81", 1) | Out-Null
$d.Content.Find.Execute("81+15=96", $true, $false, $false, $false, $false, $true, 1, $false, "3+61=64", 1) | Out-Null
$d.Content.Find.Execute("76-43=33", $true, $false, $false, $false, $false, $true, 1, $false, "41+24=65", 1) | Out-Null
$d.Content.Find.Execute("93-11=82", $true, $false, $false, $false, $false, $true, 1, $false, "91+2=93", 1) | Out-Null
$d.Content.Find.Execute("66+3=69", $true, $false, $false, $false, $false, $true, 1, $false, "37+25=62", 1) | Out-Null
$d.Content.Find.Execute("70+14=84", $true, $false, $false, $false, $false, $true, 1, $false, "88-14=74", 1) | Out-Null
$d.Content.Find.Execute("36+44=80", $true, $false, $false, $false, $false, $true, 1, $false, "38+46=84", 1) | Out-Null
$d.Content.Find.Execute("27+65=92", $true, $false, $false, $false, $false, $true, 1, $false, "90-14=76", 1) | Out-Null
$d.Content.Find.Execute("42+21=63", $true, $false, $false, $false, $false, $true, 1, $false, "82-55=27", 1) | Out-Null
$d.Content.Find.Execute("66-0=66", $true, $false, $false, $false, $false, $true, 1, $false, "82+0=82", 1) | Out-Null
$d.Content.Find.Execute("30-11=19", $true, $false, $false, $false, $false, $true, 1, $false, "93-73=20", 1) | Out-Null
$d.Content.Find.Execute("5+19=24", $true, $false, $false, $false, $false, $true, 1, $false, "95-58=37", 1) | Out-Null
$d.Content.Find.Execute("66-51=15", $true, $false, $false, $false, $false, $true, 1, $false, "75+0=75", 1) | Out-Null
$d.Content.Find.Execute("93-12=81", $true, $false, $false, $false, $false, $true, 1, $false, "27-20=7", 1) | Out-Null
$d.Content.Find.Execute("24+57=81", $true, $false, $false, $false, $false, $true, 1, $false, "26+69=95", 1) | Out-Null
$d.Content.Find.Execute("79-21=58", $true, $false, $false, $false, $false, $true, 1, $false, "55-10=45", 1) | Out-Null
$d.Content.Find.Execute("22+12=34", $true, $false, $false, $false, $false, $true, 1, $false, "48-37=11", 1) | Out-Null
$d.Content.Find.Execute("87+9=96", $true, $false, $false, $false, $false, $true, 1, $false, "96-0=96", 1) | Out-Null
$d.Content.Find.Execute("40+29=69", $true, $false, $false, $false, $false, $true, 1, $false, "56+29=85", 1) | Out-Null
$d.Content.Find.Execute("75-53=22", $true, $false, $false, $false, $false, $true, 1, $false, "56+11=67", 1) | Out-Null
$d.Content.Find.Execute("98-37=61", $true, $false, $false, $false, $false, $true, 1, $false, "1+11=12", 1) | Out-Null
$d.Content.Find.Execute("86-73=13", $true, $false, $false, $false, $false, $true, 1, $false, "84-47=37", 1) | Out-Null
$d.Content.Find.Execute("20+43=63", $true, $false, $false, $false, $false, $true, 1, $false, "31-13=18", 1) | Out-Null
$d.Content.Find.Execute("66-9=57", $true, $false, $false, $false, $false, $true, 1, $false, "57+21=78", 1) | Out-Null
$d.Content.Find.Execute("35+32=67", $true, $false, $false, $false, $false, $true, 1, $false, "45-12=33", 1) | Out-Null
$d.Content.Find.Execute("21+42=63", $true, $false, $false, $false, $false, $true, 1, $false, "36+47=83", 1) | Out-Null
$d.Content.Find.Execute("63+22=85", $true, $false, $false, $false, $false, $true, 1, $false, "68-0=68", 1) | Out-Null
$d.Content.Find.Execute("2+78=80", $true, $false, $false, $false, $false, $true, 1, $false, "86-46=40", 1) | Out-Null
$d.Content.Find.Execute("29+46=75", $true, $false, $false, $false, $false, $true, 1, $false, "89-32=57", 1) | Out-Null
$d.Content.Find.Execute("71+27=98", $true, $false, $false, $false, $false, $true, 1, $false, "7+44=51", 1) | Out-Null
$d.Content.Find.Execute("15+84=99", $true, $false, $false, $false, $false, $true, 1, $false, "52+22=74", 1) | Out-Null
$d.Content.Find.Execute("17+26=43", $true, $false, $false, $false, $false, $true, 1, $false, "76-13=63", 1) | Out-Null
$d.Content.Find.Execute("48-47=1", $true, $false, $false, $false, $false, $true, 1, $false, "94-92=2", 1) | Out-Null
$d.Content.Find.Execute("6+49=55", $true, $false, $false, $false, $false, $true, 1, $false, "99-69=30", 1) | Out-Null
$d.Content.Find.Execute("85-80=5", $true, $false, $false, $false, $false, $true, 1, $false, "15-5=10", 1) | Out-Null
$d.Content.Find.Execute("23-1=22", $true, $false, $false, $false, $false, $true, 1, $false, "8+87=95", 1) | Out-Null
$d.Content.Find.Execute("3+92=95", $true, $false, $false, $false, $false, $true, 1, $false, "39+23=62", 1) | Out-Null
$d.Content.Find.Execute("11+19=30", $true, $false, $false, $false, $false, $true, 1, $false, "79-17=62", 1) | Out-Null
$d.Content.Find.Execute("44-40=4", $true, $false, $false, $false, $false, $true, 1, $false, "16-4=12", 1) | Out-Null
